# Remove empty placeholder text frames (p:ph type="body" idx="1" with no
# visible text) that were left behind as leftover/unused frames on several
# slides. Each of these is the 2nd shape in its slide's shape tree.

$p = $ppt.ActivePresentation

# Slide number -> shape id of the empty body placeholder to delete
$targets = @{
    6  = 106
    7  = 115
    8  = 122
    9  = 130
    10 = 138
    11 = 144
    12 = 152
    13 = 160
    16 = 188
}

foreach ($slideIndex in $targets.Keys) {
    $expectedId = $targets[$slideIndex]
    $slide = $p.Slides.Item($slideIndex)

    for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.Id -eq $expectedId) {
            $shape.Delete()
            break
        }
    }
}
